$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to Text format so numeric-looking strings (e.g.
    # "0.998", "1.00") are stored as text instead of being auto-converted
    # to a number by Excel, then restore the original "Normal" style so
    # no stray number-format styling is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '92.916.43'
$ws.Range("E2").Value = '  -3.84%  '

$ws.Range("D3").Value = '3.341.21'
$ws.Range("E3").Value = '  -3.84%  '

Set-TextValue $ws.Range('D4') '0.998'
$ws.Range("E4").Value = '  -0.49%  '

Set-TextValue $ws.Range('D5') '231.33'
$ws.Range("E5").Value = '  -6.11%  '

Set-TextValue $ws.Range('D6') '619.67'
$ws.Range("E6").Value = '  -4.44%  '

Set-TextValue $ws.Range('D7') '1.37'
$ws.Range("E7").Value = '  -4.21%  '

Set-TextValue $ws.Range('D8') '0.384'
$ws.Range("E8").Value = '  -7.09%  '

Set-TextValue $ws.Range('D9') '0.999'
$ws.Range("E9").Value = '  -0.21%  '

Set-TextValue $ws.Range('D10') '0.926'
$ws.Range("E10").Value = '  -7.64%  '

$ws.Range("D11").Value = '3.342.09'
$ws.Range("E11").Value = '  -3.72%  '

Set-TextValue $ws.Range('D12') '41.93'
$ws.Range("E12").Value = '  -4.70%  '

Set-TextValue $ws.Range('D13') '0.191'
$ws.Range("E13").Value = '  -4.80%  '

Set-TextValue $ws.Range('D14') '5.99'
$ws.Range("E14").Value = '  -3.17%  '

$ws.Range("D15").Value = '92.485.78'
$ws.Range("E15").Value = '  -4.23%  '

$ws.Range("D16").Value = '3.966.08'
$ws.Range("E16").Value = '  -3.98%  '

Set-TextValue $ws.Range('D17') '0.0000242'
$ws.Range("E17").Value = '  -4.51%  '

Set-TextValue $ws.Range('D18') '7.96'
$ws.Range("E18").Value = '  -8.38%  '

$ws.Range("D19").Value = '3.337.58'
$ws.Range("E19").Value = '  -5.17%  '

Set-TextValue $ws.Range('D20') '17.27'
$ws.Range("E20").Value = '  -5.63%  '

Set-TextValue $ws.Range('D21') '11.18'
$ws.Range("E21").Value = '  -5.75%  '

Set-TextValue $ws.Range('D22') '488.98'
$ws.Range("E22").Value = '  -5.71%  '

Set-TextValue $ws.Range('D23') '3.28'
$ws.Range("E23").Value = '  -0.89%  '

Set-TextValue $ws.Range('D24') '0.447'
$ws.Range("E24").Value = '  -10.86%  '

Set-TextValue $ws.Range('D25') '0.0000182'
$ws.Range("E25").Value = '  -7.78%  '

Set-TextValue $ws.Range('D26') '6.09'
$ws.Range("E26").Value = '  -7.72%  '

Set-TextValue $ws.Range('D27') '89.34'
$ws.Range("E27").Value = '  -3.67%  '

$ws.Range("D28").Value = '3.511.53'
$ws.Range("E28").Value = '  -4.73%  '

Set-TextValue $ws.Range('D29') '11.59'
$ws.Range("E29").Value = '  -7.02%  '

Set-TextValue $ws.Range('D30') '0.999'
$ws.Range("E30").Value = '  +0.22%  '

Set-TextValue $ws.Range('D31') '11.13'
$ws.Range("E31").Value = '  -7.59%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D32') '2.65'
$ws.Range("E32").Value = '  -4.15%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D33') '0.135'
$ws.Range("E33").Value = '  -3.77%  '

Set-TextValue $ws.Range('D34') '0.997'
$ws.Range("E34").Value = '  -0.12%  '

Set-TextValue $ws.Range('D35') '0.171'
$ws.Range("E35").Value = '  -8.60%  '

Set-TextValue $ws.Range('D36') '28.28'
$ws.Range("E36").Value = '  -8.47%  '

Set-TextValue $ws.Range('D37') '0.527'
$ws.Range("E37").Value = '  -9.71%  '

Set-TextValue $ws.Range('D38') '529.50'
$ws.Range("E38").Value = '  +3.55%  '

$ws.Range("E39").Value = '  -0.03%  '

Set-TextValue $ws.Range('D40') '7.32'
$ws.Range("E40").Value = '  -7.25%  '

$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D41') '1.36'
$ws.Range("E41").Value = '  -7.92%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D42') '0.146'
$ws.Range("E42").Value = '  -5.14%  '

Set-TextValue $ws.Range('D43') '0.880'
$ws.Range("E43").Value = '  -2.83%  '

Set-TextValue $ws.Range('D44') '24.02'
$ws.Range("E44").Value = '  -0.95%  '

Set-TextValue $ws.Range('D45') '1.67'
$ws.Range("E45").Value = '  -2.78%  '

$ws.Range("E46").Value = '  -1.18%  '

Set-TextValue $ws.Range('D47') '0.0404'
$ws.Range("E47").Value = '  -4.63%  '

Set-TextValue $ws.Range('D48') '5.37'
$ws.Range("E48").Value = '  -3.89%  '

Set-TextValue $ws.Range('D49') '2.12'
$ws.Range("E49").Value = '  -4.08%  '

Set-TextValue $ws.Range('D50') '52.00'
$ws.Range("E50").Value = '  -3.50%  '

Set-TextValue $ws.Range('D51') '7.87'
$ws.Range("E51").Value = '  -6.86%  '
